$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.14'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.37%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '29.80'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.17%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.155'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.61%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05759'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.40%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.660'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.02%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.229'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '6.12%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8500'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.02%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8556'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-2.02%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1389'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.94%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07082'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.36%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03246'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '10.93%'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.24%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001536'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.91%'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0005935'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.51%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005934'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.86%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.521'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.43%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.204'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.94%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3167'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.55%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03372'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.56%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1311'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.25%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.503'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-3.10%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.1409'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.18%'
$ws.Range("B24").Value = 'CoinExToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.04114'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.60%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.07%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004151'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-7.96%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001199'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '1.72%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '4.12%'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.06%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1072'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.14%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.002469'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-4.61%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003539'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-38.40%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009946'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '5.13%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005474'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '7.27%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000749'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.03%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07094'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-20.25%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002465'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-10.89%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002098'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.03%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001998'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.03%'
